# ToDoList_Form.xlsx edit:
#  - "DB" sheet: rotate task list rows 20-24, append new rows 25-33
#    (many-team rotation entries), dimension grows to A1:I33.
#  - View-state fallout: "DB" becomes the active sheet/tab, "Team"
#    sheet's previous tabSelected/topLeftCell are cleared and its
#    lingering selection moves to F26.

$wb = $excel.ActiveWorkbook

$wsTeam = $wb.Worksheets.Item("Team")
$wsDb   = $wb.Worksheets.Item("DB")

# ---------------------------------------------------------------------
# 1. Update the existing DB rows 20-24 (values shift down / change).
# ---------------------------------------------------------------------

# Row 20
$wsDb.Cells.Item(20,1).Value2 = "06-00-00"
$wsDb.Cells.Item(20,2).Value2 = "경영층 특강 참석자 선정"
$wsDb.Cells.Item(20,3).Value2 = " "
$wsDb.Cells.Item(20,4).Value2 = " "
$wsDb.Cells.Item(20,5).Value2 = "진행"
$wsDb.Cells.Item(20,6).Value2 = " "
$wsDb.Cells.Item(20,7).Value2 = " "
$wsDb.Cells.Item(20,8).Value2 = " "
$wsDb.Cells.Item(20,9).Value2 = 1

# Row 21
$wsDb.Cells.Item(21,1).Value2 = "07-00-00"
$wsDb.Cells.Item(21,2).Value2 = "센터 주간이슈 논의회"
$wsDb.Cells.Item(21,3).Value2 = " "
$wsDb.Cells.Item(21,4).Value2 = " "
$wsDb.Cells.Item(21,5).Value2 = "진행"
$wsDb.Cells.Item(21,6).Value2 = " "
$wsDb.Cells.Item(21,7).Value2 = " "
$wsDb.Cells.Item(21,8).Value2 = " "
$wsDb.Cells.Item(21,9).Value2 = 1

# Row 22 - "07-01-00" parses as a date unless the cell is pre-formatted as text
$wsDb.Cells.Item(22,1).NumberFormat = "@"
$wsDb.Cells.Item(22,1).Value2 = "07-01-00"
$wsDb.Cells.Item(22,2).Value2 = "메일 송부"
$wsDb.Cells.Item(22,3).Value2 = " "
$wsDb.Cells.Item(22,4).Value2 = " "
$wsDb.Cells.Item(22,5).Value2 = "진행"
$wsDb.Cells.Item(22,6).Value2 = " "
$wsDb.Cells.Item(22,7).Value2 = " "
$wsDb.Cells.Item(22,8).Value2 = " "
$wsDb.Cells.Item(22,9).Value2 = 2

# Row 23 - "07-01-01" also parses as a date unless pre-formatted as text
$wsDb.Cells.Item(23,1).NumberFormat = "@"
$wsDb.Cells.Item(23,1).Value2 = "07-01-01"
$wsDb.Cells.Item(23,2).Value2 = " "
$wsDb.Cells.Item(23,3).Value2 = " "
$wsDb.Cells.Item(23,4).Value2 = " "
$wsDb.Cells.Item(23,5).Value2 = "진행"
$wsDb.Cells.Item(23,6).Value2 = " "
$wsDb.Cells.Item(23,7).Value2 = " "
$wsDb.Cells.Item(23,8).Value2 = " "
$wsDb.Cells.Item(23,9).Value2 = 3

# Row 24
$wsDb.Cells.Item(24,1).Value2 = "08-00-00"
$wsDb.Cells.Item(24,2).Value2 = "경영층 보고"
$wsDb.Cells.Item(24,3).Value2 = " "
$wsDb.Cells.Item(24,4).Value2 = " "
$wsDb.Cells.Item(24,5).Value2 = "진행"
$wsDb.Cells.Item(24,6).Value2 = " "
$wsDb.Cells.Item(24,7).Value2 = " "
$wsDb.Cells.Item(24,8).Value2 = " "
$wsDb.Cells.Item(24,9).Value2 = 1

# ---------------------------------------------------------------------
# 2. Append the new "many teams rotation" rows 25-33.
#    Column A on these rows is also a date-like "NN-NN-00" string, so
#    pre-format the whole A25:A33 block as text before writing into it.
# ---------------------------------------------------------------------

$wsDb.Range("A25:A33").NumberFormat = "@"

$newRows = @(
    @{ Row=25; A="08-01-00"; B="완료";                               C="차체설계1팀" },
    @{ Row=26; A="08-02-00"; B=" ";                                  C="차체설계2팀" },
    @{ Row=27; A="08-03-00"; B=" ";                                  C="차체설계3팀" },
    @{ Row=28; A="08-04-00"; B=" ";                                  C="외장설계1팀" },
    @{ Row=29; A="08-05-00"; B=" ";                                  C="dkkkkkkkkkkkkkkkkk" },
    @{ Row=30; A="08-06-00"; B=" ";                                  C="daaleiw12222222222222222222222" },
    @{ Row=31; A="08-07-00"; B=" ";                                  C="aksdakalskdfasf" },
    @{ Row=32; A="08-08-00"; B=" ";                                  C="1212123124k" },
    @{ Row=33; A="08-09-00"; B=" ";                                  C="외장설계2팀" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $wsDb.Cells.Item($row,1).Value2 = $r.A
    $wsDb.Cells.Item($row,2).Value2 = $r.B
    $wsDb.Cells.Item($row,3).Value2 = $r.C
    $wsDb.Cells.Item($row,4).Value2 = " "
    $wsDb.Cells.Item($row,5).Value2 = "진행"
    $wsDb.Cells.Item($row,6).Value2 = " "
    $wsDb.Cells.Item($row,7).Value2 = " "
    $wsDb.Cells.Item($row,8).Value2 = " "
    $wsDb.Cells.Item($row,9).Value2 = 2
}

# ---------------------------------------------------------------------
# 3. View-state: Team loses tabSelected/topLeftCell, its selection is
#    left at F26; DB becomes the active/selected sheet of the workbook.
# ---------------------------------------------------------------------

$wsTeam.Activate()
$wsTeam.Range("F26").Select()

$wsDb.Activate()
